# Nieuwenhuys2017.xlsx - "specify plane to all patterns"
# Fill the "Plane" column (A) down for every row of each joint-pattern
# block on the "patterns" sheet, and mark each block with the matching
# built-in Excel cell style (20/40/60% - Accent1 for Sagittal/Coronal/
# Transverse plane respectively). Also move the active tab/selection
# from "data" to "patterns".

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")
$wsPatterns = $wb.Worksheets.Item("patterns")

# Row ranges (inclusive) for each plane block on the "patterns" sheet,
# together with the plane label and the matching built-in cell style.
$blocks = @(
    @{ First = 2;  Last = 32; Plane = "Sagittal plane";  Style = "20% - Accent1" },
    @{ First = 33; Last = 40; Plane = "Coronal plane";   Style = "40% - Accent1" },
    @{ First = 41; Last = 50; Plane = "Transverse plane"; Style = "60% - Accent1" },
    @{ First = 51; Last = 75; Plane = "Sagittal plane";  Style = "20% - Accent1" },
    @{ First = 76; Last = 79; Plane = "Coronal plane";   Style = "40% - Accent1" },
    @{ First = 80; Last = 85; Plane = "Transverse plane"; Style = "60% - Accent1" }
)

foreach ($block in $blocks) {
    $addr = "A" + $block.First + ":A" + $block.Last
    $rng = $wsPatterns.Range($addr)
    $rng.Value = $block.Plane
    $rng.Style = $block.Style
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# The edit also shifts the active tab from "data" to "patterns", scrolls
# down within "patterns" and updates its selection.
$wsPatterns.Activate()
$wsPatterns.Range("A80:A85").Select()
